$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# N1 header: remove trailing space -> "Correction"
$ws.Range("N1").Value = "Correction"

# N2:N12 were empty inline strings; set them to the text "nan" to match
# the rest of the column's placeholder values.
$ws.Range("N2:N12").Value = "nan"
